$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Recipe Types"
$ws2 = $wb.Worksheets.Item(2)   # "Recipes "

# Add the new "Users" worksheet after the last existing sheet ("Recipes ")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Users"

# Match column A width/format to the "Recipe Types" sheet
$newSheet.Columns.Item(1).ColumnWidth = $ws1.Columns.Item(1).ColumnWidth

# Copy the header row formats (bold red "Collection: ..." row, bold Field/Datatype row)
# from the "Recipe Types" sheet so we reuse the same cell styles (no new styles added).
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)
$ws1.Range("A2:B2").Copy()
$newSheet.Range("A2:B2").PasteSpecial(-4122)

# Fill in the field/datatype header row first so new shared strings are interned
# in the same order as the target workbook.
$newSheet.Range("A2").Value = "Field"
$newSheet.Range("B2").Value = "Datatype"

$newSheet.Range("A3").Value = "first_name"
$newSheet.Range("B3").Value = "String"
$newSheet.Range("A4").Value = "last_name"
$newSheet.Range("B4").Value = "String"
$newSheet.Range("A5").Value = "username"
$newSheet.Range("B5").Value = "String"
$newSheet.Range("A6").Value = "user_type"
$newSheet.Range("B6").Value = "String"

# Correct the "Recipe Types" sheet's title text (it mistakenly said "Collection: recipes")
$ws1.Range("A1").Value = "Collection: recipe_types"

# Title for the new "Users" sheet
$newSheet.Range("A1").Value = "Collection: users"

# Restore/update the selections on the existing sheets, and leave "Recipe Types"
# as the active tab/sheet.
$ws2.Activate() | Out-Null
$ws2.Range("B3").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("A4").Select() | Out-Null
